# Edit described by the diff:
#  1. Paragraph "Question" (pStyle FirstParagraph) loses its run-level
#     formatting and is split into two runs: "Questio" + "n", both with
#     an explicit-but-empty <w:rPr/>.
#  2. The "FirstParagraph" paragraph style gains an explicit font color
#     (00CC00) in its style-level <w:rPr/>.

$d = $word.ActiveDocument

# --- 1. Split "Question" into "Questio" + "n", stripping formatting ---

$rng = $d.Content
$rng.Find.Execute("Question", $true, $false, $false, $false, $false, $true, `
                   1, $false, "", 0)
$qStart = $rng.Start

# Remove the existing richly-formatted run entirely.
$rng.Delete()

# Re-insert the plain text at the same spot; a fresh insertion point has
# no run formatting of its own.
$whole = $d.Range($qStart, $qStart)
$whole.InsertAfter("Question")

# Touch the "Questio" sub-range's Font so the engine materializes an
# explicit (empty) <w:rPr/> on its run instead of omitting rPr entirely.
$partA = $d.Range($qStart, $qStart + 7)
$partA.Font.Bold = $true
$partA.Font.Bold = $false

# Do the same for the trailing "n" sub-range; touching its formatting
# also keeps it from being re-merged into the previous run.
$partB = $d.Range($qStart + 7, $qStart + 8)
$partB.Font.Bold = $true
$partB.Font.Bold = $false

# --- 2. Give the "FirstParagraph" style an explicit green font color ---

$style = $d.Styles("FirstParagraph")
$style.Font.Color = 52224
